# Tasks reader test fixture update — "Added ability to skip row in excel reader"
#
# Net effect on the "tasks" sheet:
#   - rows for TASK-3 / TASK-4 shift down by one row (row 6 becomes free / skipped)
#   - a new centered section header row "Some other tasks" (merged A:G) is inserted
#   - a new blank centered row follows it
#   - one row is left completely empty (the "skipped" row) before TASK-5 / TASK-6
#   - selection moves to the now-empty A11
#   - a new shared string "Some other tasks" is introduced

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tasks")

# Shift old row 6 (TASK-3) .. row 9 (TASK-6) down by one, opening a gap at row 6.
$ws.Rows("6").Insert() | Out-Null

# Open three more rows right before the old TASK-5 row (now row 9), for the new
# section header, a blank spacer row, and one fully-skipped row.
$ws.Rows("9:11").Insert() | Out-Null

# Re-apply the centered header formatting to the title rows (A1:F2) so it keeps
# matching the freshly introduced section-header style.
$ws.Range("A1:F2").HorizontalAlignment = -4108

# New merged, centered section header row.
$ws.Range("A9:G9").HorizontalAlignment = -4108
$ws.Range("A9:G9").Merge() | Out-Null
$ws.Range("A9").Value = "Some other tasks"

# New blank centered spacer row underneath the header.
$ws.Range("A10:G10").HorizontalAlignment = -4108

# Row 11 is intentionally left empty (the "skipped" row).

# Park the selection on the now-empty row, matching where the author left off.
$ws.Range("A11").Select() | Out-Null
